$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily attendance processing - 2025-10-18 10:18:32
# Applies 74 cell updates: "Recorded By" (col G) name-order normalization
# plus a handful of corrected attendance counts (col H) and recomputed
# average-attendance percentages (col L / col S).

$updates = @{
    'G2' = 'system, System, backup@backdoor.com'
    'G3' = 'System, dnasr281@gmail.com'
    'G6' = 'System, dnasr281@gmail.com'
    'G10' = 'System, dnasr281@gmail.com'
    'L10' = '72.5%'
    'G11' = 'System, dnasr281@gmail.com'
    'G12' = 'System, dnasr281@gmail.com'
    'G13' = 'System, dnasr281@gmail.com'
    'G14' = 'System, dnasr281@gmail.com'
    'G15' = 'System, dnasr281@gmail.com'
    'S15' = '70.1%'
    'G17' = 'System, dnasr281@gmail.com'
    'G18' = 'System, dnasr281@gmail.com'
    'H18' = '38/53'
    'G19' = 'System, dnasr281@gmail.com'
    'G20' = 'System, dnasr281@gmail.com'
    'S20' = '83.2%'
    'G29' = 'system, System, backup@backdoor.com'
    'G30' = 'System, dnasr281@gmail.com'
    'G33' = 'System, dnasr281@gmail.com'
    'G37' = 'System, dnasr281@gmail.com'
    'G38' = 'System, dnasr281@gmail.com'
    'G39' = 'System, dnasr281@gmail.com'
    'G40' = 'System, dnasr281@gmail.com'
    'G41' = 'System, dnasr281@gmail.com'
    'G42' = 'System, dnasr281@gmail.com'
    'G44' = 'System, dnasr281@gmail.com'
    'G45' = 'System, dnasr281@gmail.com'
    'G46' = 'System, dnasr281@gmail.com'
    'G47' = 'System, dnasr281@gmail.com'
    'G56' = 'system, System, backup@backdoor.com'
    'G57' = 'System, dnasr281@gmail.com'
    'G60' = 'System, dnasr281@gmail.com'
    'G64' = 'System, dnasr281@gmail.com'
    'G65' = 'System, dnasr281@gmail.com'
    'G66' = 'System, dnasr281@gmail.com'
    'G67' = 'System, dnasr281@gmail.com'
    'G68' = 'System, dnasr281@gmail.com'
    'G69' = 'System, dnasr281@gmail.com'
    'G71' = 'System, dnasr281@gmail.com'
    'G72' = 'System, dnasr281@gmail.com'
    'G73' = 'System, dnasr281@gmail.com'
    'G74' = 'System, dnasr281@gmail.com'
    'G86' = 'System, dnasr281@gmail.com'
    'G87' = 'System, dnasr281@gmail.com'
    'G88' = 'System, dnasr281@gmail.com'
    'G89' = 'System, dnasr281@gmail.com'
    'G90' = 'admin@admin.com, dnasr281@gmail.com'
    'G93' = 'System, dnasr281@gmail.com'
    'G95' = 'System, dnasr281@gmail.com'
    'G96' = 'System, dnasr281@gmail.com'
    'G97' = 'System, dnasr281@gmail.com'
    'G99' = 'System, dnasr281@gmail.com'
    'G112' = 'System, dnasr281@gmail.com'
    'G113' = 'System, dnasr281@gmail.com'
    'G114' = 'System, dnasr281@gmail.com'
    'G115' = 'System, dnasr281@gmail.com'
    'G116' = 'admin@admin.com, dnasr281@gmail.com'
    'G119' = 'System, dnasr281@gmail.com'
    'G121' = 'System, dnasr281@gmail.com'
    'G122' = 'System, dnasr281@gmail.com'
    'G123' = 'System, dnasr281@gmail.com'
    'G125' = 'System, dnasr281@gmail.com'
    'H135' = '57/57'
    'G138' = 'System, dnasr281@gmail.com'
    'G139' = 'System, dnasr281@gmail.com'
    'G140' = 'System, dnasr281@gmail.com'
    'G141' = 'System, dnasr281@gmail.com'
    'G142' = 'admin@admin.com, dnasr281@gmail.com'
    'G145' = 'System, dnasr281@gmail.com'
    'G147' = 'System, dnasr281@gmail.com'
    'G148' = 'System, dnasr281@gmail.com'
    'G149' = 'System, dnasr281@gmail.com'
    'G151' = 'System, dnasr281@gmail.com'
}

# Cells holding a percentage written as literal text (e.g. "72.5%") must be
# forced to the Text number format first, otherwise Excel auto-converts the
# string into a numeric percentage (0.725) on assignment.
$textCells = @('L10', 'S15', 'S20')

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($textCells -contains $cellRef) {
        $range.NumberFormat = "@"
    }
    $range.Value = $updates[$cellRef]
}

